# Strings.xlsx update:
#  - add "save title" / "S A V E" and "load title" / "L O A D" rows just
#    before the existing "save" row on the localization sheet
#  - add a "challenge warning" / "Warning: some levels may require advanced
#    techniques." row just after "challenge title" / "C H A L L E N G E"
#
# These are simple row insertions - every row below each insertion point
# shifts down and keeps its original content untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- insert "save title" / "load title" rows above the "save" row ---------
# Before the edit, row 75 holds the "save" key ("save" / "Save").
$ws.Rows(75).Insert()
$ws.Cells.Item(75, 1).Value2 = "save title"
$ws.Cells.Item(75, 2).Value2 = "S A V E"

$ws.Rows(76).Insert()
$ws.Cells.Item(76, 1).Value2 = "load title"
$ws.Cells.Item(76, 2).Value2 = "L O A D"

# --- insert "challenge warning" row below "challenge title" ---------------
# After the two rows above were inserted, the old "challenge title" row
# (originally row 86) is now at row 88, so the new row goes in at 89.
$ws.Rows(89).Insert()
$ws.Cells.Item(89, 1).Value2 = "challenge warning"
$ws.Cells.Item(89, 2).Value2 = "Warning: some levels may require advanced techniques."

# Keep the sheet's selection on the same logical row as the author left it.
[void]$ws.Range("A77").Select()
